$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextCell($cell, $value) {
    # Force text entry so numeric-looking strings (e.g. "1.005",
    # "225.56") are kept verbatim instead of being parsed into
    # floating point numbers, then drop back to the default/
    # unstyled cell format (matches source: no explicit style).
    $cell.NumberFormat = "@"
    $cell.Value = $value
    $cell.Style = "Normal"
}

Set-TextCell $ws.Range('D2') '27.517.91'
Set-TextCell $ws.Range('E2') '  +5.42%  '
Set-TextCell $ws.Range('D3') '1.726.08'
Set-TextCell $ws.Range('E3') '  +4.49%  '
Set-TextCell $ws.Range('D4') '1.005'
Set-TextCell $ws.Range('E4') '  +0.19%  '
Set-TextCell $ws.Range('D5') '225.56'
Set-TextCell $ws.Range('E5') '  +3.30%  '
Set-TextCell $ws.Range('D6') '0.5356'
Set-TextCell $ws.Range('E6') '  +2.88%  '
Set-TextCell $ws.Range('D7') '1.005'
Set-TextCell $ws.Range('E7') '  +0.12%  '
Set-TextCell $ws.Range('D8') '0.2664'
Set-TextCell $ws.Range('E8') '  +0.74%  '
Set-TextCell $ws.Range('D9') '0.06595'
Set-TextCell $ws.Range('E9') '  +4.10%  '
Set-TextCell $ws.Range('D10') '21.70'
Set-TextCell $ws.Range('E10') '  +6.55%  '
Set-TextCell $ws.Range('D11') '0.07716'
Set-TextCell $ws.Range('E12') '  -0.17%  '
Set-TextCell $ws.Range('D13') '1.739.92'
Set-TextCell $ws.Range('E13') '  +4.68%  '
Set-TextCell $ws.Range('D14') '1.964.21'
Set-TextCell $ws.Range('E14') '  +4.51%  '
Set-TextCell $ws.Range('D15') '0.5842'
Set-TextCell $ws.Range('E15') '  +4.40%  '
Set-TextCell $ws.Range('D16') '0.0₅8292'
Set-TextCell $ws.Range('E16') '  +1.67%  '
Set-TextCell $ws.Range('D17') '67.95'
Set-TextCell $ws.Range('E17') '  +3.95%  '
Set-TextCell $ws.Range('D18') '27.545.11'
Set-TextCell $ws.Range('E18') '  +5.50%  '
Set-TextCell $ws.Range('D19') '219.50'
Set-TextCell $ws.Range('E19') '  +15.10%  '
Set-TextCell $ws.Range('D20') '1.005'
Set-TextCell $ws.Range('E20') '  +0.10%  '
Set-TextCell $ws.Range('D21') '4.727'
Set-TextCell $ws.Range('E21') '  +2.10%  '
Set-TextCell $ws.Range('E22') '  +1.39%  '
Set-TextCell $ws.Range('D23') '6.093'
Set-TextCell $ws.Range('E23') '  +2.64%  '
Set-TextCell $ws.Range('E24') '  +0.14%  '
Set-TextCell $ws.Range('D25') '148.38'
Set-TextCell $ws.Range('E25') '  +2.65%  '
Set-TextCell $ws.Range('D26') '1.715'
Set-TextCell $ws.Range('E26') '  +14.04%  '
Set-TextCell $ws.Range('D27') '0.1235'
Set-TextCell $ws.Range('E27') '  +3.76%  '
Set-TextCell $ws.Range('D28') '7.410'
Set-TextCell $ws.Range('E28') '  +2.68%  '
Set-TextCell $ws.Range('D29') '16.66'
Set-TextCell $ws.Range('E29') '  +4.60%  '
Set-TextCell $ws.Range('D30') '0.05562'
Set-TextCell $ws.Range('E30') '  +1.44%  '
Set-TextCell $ws.Range('E31') '  +2.51%  '
Set-TextCell $ws.Range('D32') '3.550'
Set-TextCell $ws.Range('E32') '  +2.94%  '
Set-TextCell $ws.Range('D33') '3.457'
Set-TextCell $ws.Range('E33') '  +2.89%  '
Set-TextCell $ws.Range('D34') '1.659'
Set-TextCell $ws.Range('E34') '  +6.38%  '
Set-TextCell $ws.Range('D35') '0.9609'
Set-TextCell $ws.Range('E35') '  +1.26%  '
Set-TextCell $ws.Range('E36') '  +1.50%  '
Set-TextCell $ws.Range('D37') '2.432'
Set-TextCell $ws.Range('E37') '  +1.48%  '
Set-TextCell $ws.Range('D38') '0.5958'
Set-TextCell $ws.Range('E38') '  +5.62%  '
Set-TextCell $ws.Range('D39') '0.01653'
Set-TextCell $ws.Range('E39') '  +4.80%  '
Set-TextCell $ws.Range('D40') '5.932'
Set-TextCell $ws.Range('E40') '  +1.17%  '
Set-TextCell $ws.Range('B41') 'Maker'
Set-TextCell $ws.Range('C41') 'https://coinranking.com/coin/qFakph2rpuMOL+maker-mkr'
Set-TextCell $ws.Range('D41') '1.057.15'
Set-TextCell $ws.Range('E41') '  +2.64%  '
Set-TextCell $ws.Range('B42') 'TrustWalletToken'
Set-TextCell $ws.Range('C42') 'https://coinranking.com/coin/Hm3OlynlC+trustwallettoken-twt'
Set-TextCell $ws.Range('D42') '0.8547'
Set-TextCell $ws.Range('E42') '  +2.63%  '
Set-TextCell $ws.Range('E43') '  +0.12%  '
Set-TextCell $ws.Range('D44') '101.49'
Set-TextCell $ws.Range('E44') '  +0.61%  '
Set-TextCell $ws.Range('D45') '1.870.41'
Set-TextCell $ws.Range('E45') '  +4.37%  '
Set-TextCell $ws.Range('E46') '  +6.55%  '
Set-TextCell $ws.Range('D47') '58.98'
Set-TextCell $ws.Range('E47') '  +2.22%  '
Set-TextCell $ws.Range('D48') '8.218'
Set-TextCell $ws.Range('E48') '  +2.67%  '
Set-TextCell $ws.Range('D49') '0.4440'
Set-TextCell $ws.Range('E49') '  +2.33%  '
Set-TextCell $ws.Range('D50') '1.002'
Set-TextCell $ws.Range('E50') '  +0.16%  '
Set-TextCell $ws.Range('D51') '0.06558'
Set-TextCell $ws.Range('E51') '  +11.64%  '
